$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 and row 3: Fecha (D), Volumen (M),
# Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P),
# Precio $/Kg (S)

$row2_D = $ws.Range("D2").Value2
$row2_M = $ws.Range("M2").Value2
$row2_N = $ws.Range("N2").Value2
$row2_O = $ws.Range("O2").Value2
$row2_P = $ws.Range("P2").Value2
$row2_S = $ws.Range("S2").Value2

$row3_D = $ws.Range("D3").Value2
$row3_M = $ws.Range("M3").Value2
$row3_N = $ws.Range("N3").Value2
$row3_O = $ws.Range("O3").Value2
$row3_P = $ws.Range("P3").Value2
$row3_S = $ws.Range("S3").Value2

$ws.Range("D2").Value = $row3_D
$ws.Range("M2").Value = $row3_M
$ws.Range("N2").Value = $row3_N
$ws.Range("O2").Value = $row3_O
$ws.Range("P2").Value = $row3_P
$ws.Range("S2").Value = $row3_S

$ws.Range("D3").Value = $row2_D
$ws.Range("M3").Value = $row2_M
$ws.Range("N3").Value = $row2_N
$ws.Range("O3").Value = $row2_O
$ws.Range("P3").Value = $row2_P
$ws.Range("S3").Value = $row2_S
